$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

    $ws.Range("B2").Value = 14.48390181660633
    $ws.Range("C2").Value = 5.734302902197689
    $ws.Range("D2").Value = 7.290883180393839
    $ws.Range("F2").Value = 42.15903205432626
    $ws.Range("G2").Value = 50.20208621533102
    $ws.Range("H2").Value = 19.70584663822526
    $ws.Range("I2").Value = 31.28105683502368
    $ws.Range("J2").Value = 10.54022212088717
    $ws.Range("K2").Value = 11.90429623765264
    $ws.Range("L2").Value = 11.80628642405434
    $ws.Range("B3").Value = 14.31168369448204
    $ws.Range("C3").Value = 5.633383999484288
    $ws.Range("D3").Value = 7.280207899678215
    $ws.Range("F3").Value = 42.2157933931824
    $ws.Range("G3").Value = 50.25317137566486
    $ws.Range("H3").Value = 19.75201204293856
    $ws.Range("I3").Value = 31.35940082269659
    $ws.Range("J3").Value = 10.5602361602346
    $ws.Range("K3").Value = 11.78577553237878
    $ws.Range("L3").Value = 11.8085745886991
    $ws.Range("B4").Value = 14.20793119382198
    $ws.Range("C4").Value = 5.5695452749381
    $ws.Range("D4").Value = 7.274563744765454
    $ws.Range("F4").Value = 42.25886651814269
    $ws.Range("G4").Value = 50.29609691813972
    $ws.Range("H4").Value = 19.78326307820607
    $ws.Range("I4").Value = 31.41259177609564
    $ws.Range("J4").Value = 10.57336073849569
    $ws.Range("K4").Value = 11.7147955213377
    $ws.Range("L4").Value = 11.81162098893006
    $ws.Range("B5").Value = 14.1662002709966
    $ws.Range("C5").Value = 5.543074460698623
    $ws.Range("D5").Value = 7.27249455504572
    $ws.Range("F5").Value = 42.27848402540059
    $ws.Range("G5").Value = 50.31649062120639
    $ws.Range("H5").Value = 19.79672831840612
    $ws.Range("I5").Value = 31.43554541591072
    $ws.Range("J5").Value = 10.57891974718324
    $ws.Range("K5").Value = 11.68635161601517
    $ws.Range("L5").Value = 11.81327618925224
    $ws.Range("B6").Value = 14.15930544742839
    $ws.Range("C6").Value = 5.538651918980646
    $ws.Range("D6").Value = 7.272164963419622
    $ws.Range("F6").Value = 42.2818661265383
    $ws.Range("G6").Value = 50.32005202389072
    $ws.Range("H6").Value = 19.79900830151351
    $ws.Range("I6").Value = 31.43943399021359
    $ws.Range("J6").Value = 10.57985555147103
    $ws.Range("K6").Value = 11.68165844998615
    $ws.Range("L6").Value = 11.81357605175665
    $ws.Range("B7").Value = 14.20736610932672
    $ws.Range("C7").Value = 5.569190106014324
    $ws.Range("D7").Value = 7.274534901826378
    $ws.Range("F7").Value = 42.25912272990677
    $ws.Range("G7").Value = 50.29636021601369
    $ws.Range("H7").Value = 19.78344171916908
    $ws.Range("I7").Value = 31.41289616467514
    $ws.Range("J7").Value = 10.57343485577539
    $ws.Range("K7").Value = 11.71440992976992
    $ws.Range("L7").Value = 11.81164163496183
    $ws.Range("B8").Value = 14.424137402614
    $ws.Range("C8").Value = 5.699903959232902
    $ws.Range("D8").Value = 7.287014536531511
    $ws.Range("F8").Value = 42.17689566400963
    $ws.Range("G8").Value = 50.21729837579961
    $ws.Range("H8").Value = 19.72116126444629
    $ws.Range("I8").Value = 31.30701313758853
    $ws.Range("J8").Value = 10.54694974500165
    $ws.Range("K8").Value = 11.8630768059803
    $ws.Range("L8").Value = 11.8067353362916
    $ws.Range("B9").Value = 14.8628426492321
    $ws.Range("C9").Value = 5.940731989575024
    $ws.Range("D9").Value = 7.31862420913654
    $ws.Range("F9").Value = 42.08097430334333
    $ws.Range("G9").Value = 50.15417050321842
    $ws.Range("H9").Value = 19.62209757175284
    $ws.Range("I9").Value = 31.13980320828593
    $ws.Range("J9").Value = 10.50162471425708
    $ws.Range("K9").Value = 12.16746092488493
    $ws.Range("L9").Value = 11.81009045342446
    $ws.Range("B10").Value = 15.19052409573797
    $ws.Range("C10").Value = 6.107428659658203
    $ws.Range("D10").Value = 7.346075743123828
    $ws.Range("F10").Value = 42.05042303039014
    $ws.Range("G10").Value = 50.16402248584427
    $ws.Range("H10").Value = 19.56339922929989
    $ws.Range("I10").Value = 31.04167819198615
    $ws.Range("J10").Value = 10.47232770797534
    $ws.Range("K10").Value = 12.39704992098751
    $ws.Range("L10").Value = 11.82039668327497
    $ws.Range("B11").Value = 15.3401378413108
    $ws.Range("C11").Value = 6.180888058898946
    $ws.Range("D11").Value = 7.359453642244418
    $ws.Range("F11").Value = 42.04520194441431
    $ws.Range("G11").Value = 50.18073103711689
    $ws.Range("H11").Value = 19.53975780213139
    $ws.Range("I11").Value = 31.00242260921021
    $ws.Range("J11").Value = 10.45986317852112
    $ws.Range("K11").Value = 12.50238536527582
    $ws.Range("L11").Value = 11.82677062626928
    $ws.Range("B12").Value = 15.39682050077122
    $ws.Range("C12").Value = 6.208352418012448
    $ws.Range("D12").Value = 7.364644872852355
    $ws.Range("F12").Value = 42.04447207968403
    $ws.Range("G12").Value = 50.18881522344854
    $ws.Range("H12").Value = 19.53124579714212
    $ws.Range("I12").Value = 30.98833259540141
    $ws.Range("J12").Value = 10.45526681201497
    $ws.Range("K12").Value = 12.5423678919961
    $ws.Range("L12").Value = 11.82942495199262
    $ws.Range("B13").Value = 15.3846124821165
    $ws.Range("C13").Value = 6.202453365821094
    $ws.Range("D13").Value = 7.363521319053384
    $ws.Range("F13").Value = 42.04457381214259
    $ws.Range("G13").Value = 50.18699603910205
    $ws.Range("H13").Value = 19.53305941355794
    $ws.Range("I13").Value = 30.99133263408098
    $ws.Range("J13").Value = 10.45625122728877
    $ws.Range("K13").Value = 12.53375328937156
    $ws.Range("L13").Value = 11.82884262139518
    $ws.Range("B14").Value = 15.34480090661089
    $ws.Range("C14").Value = 6.183154711651023
    $ws.Range("D14").Value = 7.35987823223157
    $ws.Range("F14").Value = 42.04511690641331
    $ws.Range("G14").Value = 50.18136092645781
    $ws.Range("H14").Value = 19.53904868438324
    $ws.Range("I14").Value = 31.00124787157685
    $ws.Range("J14").Value = 10.45948255618379
    $ws.Range("K14").Value = 12.5056730577477
    $ws.Range("L14").Value = 11.82698418661618
    $ws.Range("B15").Value = 15.320417184316
    $ws.Range("C15").Value = 6.171287397228877
    $ws.Range("D15").Value = 7.357662974126425
    $ws.Range("F15").Value = 42.04561197022579
    $ws.Range("G15").Value = 50.1781380140919
    $ws.Range("H15").Value = 19.54277466041811
    $ws.Range("I15").Value = 31.00742223467044
    $ws.Range("J15").Value = 10.46147793365497
    $ws.Range("K15").Value = 12.48848436987304
    $ws.Range("L15").Value = 11.82587712794731
    $ws.Range("B16").Value = 15.18075309209901
    $ws.Range("C16").Value = 6.102579140085025
    $ws.Range("D16").Value = 7.345219149733165
    $ws.Range("F16").Value = 42.05093881587519
    $ws.Range("G16").Value = 50.16317650877823
    $ws.Range("H16").Value = 19.56500587390627
    $ws.Range("I16").Value = 31.04435204940928
    $ws.Range("J16").Value = 10.47315962287213
    $ws.Range("K16").Value = 12.39018105729157
    $ws.Range("L16").Value = 11.82001390423551
    $ws.Range("B17").Value = 15.095177303555
    $ws.Range("C17").Value = 6.059812552906137
    $ws.Range("D17").Value = 7.33781136104872
    $ws.Range("F17").Value = 42.05642888774096
    $ws.Range("G17").Value = 50.15712929865843
    $ws.Range("H17").Value = 19.57942820331455
    $ws.Range("I17").Value = 31.06838670700216
    $ws.Range("J17").Value = 10.48054667044314
    $ws.Range("K17").Value = 12.33007970598896
    $ws.Range("L17").Value = 11.81684753388968
    $ws.Range("B18").Value = 15.04601117615484
    $ws.Range("C18").Value = 6.034991774064684
    $ws.Range("D18").Value = 7.333634540885697
    $ws.Range("F18").Value = 42.06040340770205
    $ws.Range("G18").Value = 50.15480211015307
    $ws.Range("H18").Value = 19.58801161867721
    $ws.Range("I18").Value = 31.0827173121496
    $ws.Range("J18").Value = 10.4848767413291
    $ws.Range("K18").Value = 12.29559721212835
    $ws.Range("L18").Value = 11.81518508890723
    $ws.Range("B19").Value = 15.029375379095
    $ws.Range("C19").Value = 6.026550038883015
    $ws.Range("D19").Value = 7.33223483860647
    $ws.Range("F19").Value = 42.06188939733892
    $ws.Range("G19").Value = 50.15421186381186
    $ws.Range("H19").Value = 19.59096728076996
    $ws.Range("I19").Value = 31.08765636039317
    $ws.Range("J19").Value = 10.48635679341201
    $ws.Range("K19").Value = 12.28393788932855
    $ws.Range("L19").Value = 11.81464953274648
    $ws.Range("B20").Value = 15.10428169938308
    $ws.Range("C20").Value = 6.064388254078749
    $ws.Range("D20").Value = 7.338591263979809
    $ws.Range("F20").Value = 42.05575992822619
    $ws.Range("G20").Value = 50.15765391366929
    $ws.Range("H20").Value = 19.57786310448741
    $ws.Range("I20").Value = 31.0657757454357
    $ws.Range("J20").Value = 10.4797519018025
    $ws.Range("K20").Value = 12.33646892407827
    $ws.Range("L20").Value = 11.81716817972352
    $ws.Range("B21").Value = 15.35649420210394
    $ws.Range("C21").Value = 6.188832870065345
    $ws.Range("D21").Value = 7.360944916138407
    $ws.Range("F21").Value = 42.04492354363244
    $ws.Range("G21").Value = 50.18296842700011
    $ws.Range("H21").Value = 19.53727753294369
    $ws.Range("I21").Value = 30.9983144767742
    $ws.Range("J21").Value = 10.45853008302621
    $ws.Range("K21").Value = 12.51391861564983
    $ws.Range("L21").Value = 11.8275235372015
    $ws.Range("B22").Value = 15.52146479330035
    $ws.Range("C22").Value = 6.268101314780925
    $ws.Range("D22").Value = 7.376283407200594
    $ws.Range("F22").Value = 42.04511055422696
    $ws.Range("G22").Value = 50.20975264185844
    $ws.Range("H22").Value = 19.51332016432832
    $ws.Range("I22").Value = 30.95874376321044
    $ws.Range("J22").Value = 10.4453811211224
    $ws.Range("K22").Value = 12.63042612099064
    $ws.Range("L22").Value = 11.83569316881206
    $ws.Range("B23").Value = 15.43342138278788
    $ws.Range("C23").Value = 6.225986830648924
    $ws.Range("D23").Value = 7.368031167254252
    $ws.Range("F23").Value = 42.04434594021772
    $ws.Range("G23").Value = 50.1945212222631
    $ws.Range("H23").Value = 19.52587161627825
    $ws.Range("I23").Value = 30.97944949451664
    $ws.Range("J23").Value = 10.45233315224623
    $ws.Range("K23").Value = 12.56820627143763
    $ws.Range("L23").Value = 11.83120522367231
    $ws.Range("B24").Value = 15.10016549671892
    $ws.Range("C24").Value = 6.06232030661165
    $ws.Range("D24").Value = 7.338238414207732
    $ws.Range("F24").Value = 42.05605981608517
    $ws.Range("G24").Value = 50.15741315466157
    $ws.Range("H24").Value = 19.5785697766382
    $ws.Range("I24").Value = 31.06695456399237
    $ws.Range("J24").Value = 10.48011095764101
    $ws.Range("K24").Value = 12.33358013612585
    $ws.Range("L24").Value = 11.81702272363859
    $ws.Range("B25").Value = 14.74300647259897
    $ws.Range("C25").Value = 5.877325638134556
    $ws.Range("D25").Value = 7.309321107346801
    $ws.Range("F25").Value = 42.0999174773525
    $ws.Range("G25").Value = 50.16138266557103
    $ws.Range("H25").Value = 19.6464250244593
    $ws.Range("I25").Value = 31.18070149434047
    $ws.Range("J25").Value = 10.5131813295657
    $ws.Range("K25").Value = 12.08393356055824
    $ws.Range("L25").Value = 11.80780026918106

